{"js": "// Fill in the \"Completed By\" / \"Time Taken\" cells for the \"Mine generation\"\n// and \"Adjacent mine counter\" rows in the \"User Story 2 \u2014 Initial Minefield\n// Generation\" task table. Both rows get \"Skylar Franz\" as the author; the\n// first logs 65 minutes, the second 85 minutes.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The task-logging table under User Story 2 is the second table in the\n// document (index 1): Tasks | Completed By | Time Taken.\nconst table = tables.items[1];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst edits = [\n  { taskText: \"Mine generation\", name: \"Skylar Franz\", time: \"65 minutes\" },\n  { taskText: \"Adjacent mine counter\", name: \"Skylar Franz\", time: \"85 minutes\" },\n];\n\n// Load cell + body text for every row so we can match rows by their task\n// name rather than assuming a fixed index.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nfor (const edit of edits) {\n  const row = rows.items.find(\n    (r) => r.cells.items[0].body.text.trim() === edit.taskText\n  );\n  if (!row) {\n    throw new Error(`Could not find row for task \"${edit.taskText}\"`);\n  }\n\n  const nameCell = row.cells.items[1];\n  const timeCell = row.cells.items[2];\n\n  // insertText(..., End) on the (empty) cell paragraph appends a run while\n  // preserving the paragraph's existing formatting (widowControl/spacing/\n  // bold paragraph mark); then bolding the inserted range produces a run\n  // with <w:rPr><w:b/></w:rPr>, matching the rest of the table's rows.\n  const nameRange = nameCell.body.insertText(edit.name, Word.InsertLocation.end);\n  nameRange.font.bold = true;\n\n  const timeRange = timeCell.body.insertText(edit.time, Word.InsertLocation.end);\n  timeRange.font.bold = true;\n}\n\nawait context.sync();\n", "ps1": "# Fill in the \"Completed By\" / \"Time Taken\" cells for the \"Mine generation\"\n# and \"Adjacent mine counter\" rows in the \"User Story 2 \u2014 Initial Minefield\n# Generation\" task table. Both rows get \"Skylar Franz\" as the author; the\n# first logs 65 minutes, the second 85 minutes.\n\nfunction Clean-CellText($text) {\n  # Word COM cell Range.Text is terminated with a paragraph mark (CR, 13)\n  # and/or cell mark (BEL, 7); strip those plus normal whitespace so we can\n  # compare against plain task-name strings.\n  return $text.TrimEnd([char]13, [char]7).Trim()\n}\n\n$d = $word.ActiveDocument\n\n$edits = @{\n  \"Mine generation\"       = @{ Name = \"Skylar Franz\"; Time = \"65 minutes\" };\n  \"Adjacent mine counter\" = @{ Name = \"Skylar Franz\"; Time = \"85 minutes\" };\n}\n\nforeach ($tbl in $d.Tables) {\n  for ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $row = $tbl.Rows.Item($r)\n    if ($row.Cells.Count -lt 3) { continue }\n\n    $taskText = Clean-CellText $row.Cells.Item(1).Range.Text\n    if ($edits.ContainsKey($taskText)) {\n      $info = $edits[$taskText]\n\n      # Setting Range.Text on the (empty) cell keeps the cell's existing\n      # paragraph formatting (widowControl/spacing/bold paragraph mark);\n      # then bolding the range gives the inserted run <w:rPr><w:b/></w:rPr>,\n      # matching the rest of the table's rows.\n      $nameCell = $row.Cells.Item(2)\n      $nameCell.Range.Text = $info.Name\n      $nameCell.Range.Font.Bold = 1\n\n      $timeCell = $row.Cells.Item(3)\n      $timeCell.Range.Text = $info.Time\n      $timeCell.Range.Font.Bold = 1\n    }\n  }\n}\n"}
